$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.25555682182312
$ws.Range("B1").Value = 2.334315776824951
$ws.Range("C1").Value = 3.398365020751953
$ws.Range("D1").Value = 3.842350006103516
$ws.Range("E1").Value = 1.073539257049561
